$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$priceUpdates = @{
    3 = 8.710000000000001
    4 = 6.51
    5 = 0.72
    6 = 6.09
    7 = 0.26
    8 = 0.31
    9 = 2.06
    11 = 0.07000000000000001
    13 = 0.17
    14 = 0.1
    16 = 1.26
    17 = 0.5600000000000001
    18 = 0.09
    21 = 0.34
    22 = 0.25
    24 = 0.17
    25 = 10.38
    26 = 0.22
    28 = 0.48
    29 = 6.44
    30 = 0.2
    31 = 5.7
    32 = 0.17
    33 = 0.2
    34 = 0.5
    35 = 0.54
    36 = 0.32
    37 = 3.35
    39 = 0.33
    40 = 2.24
    41 = 5.67
    42 = 0.38
    43 = 4.13
    44 = 20.41
    45 = 0.32
    46 = 0.18
    47 = 1
    48 = 0.25
    49 = 0.16
    50 = 0.11
    51 = 0.52
    52 = 0.67
    53 = 0.62
    55 = 5.6
    56 = 0.62
    57 = 0.34
    59 = 0.28
    60 = 0.32
    62 = 0.13
    63 = 0.44
    65 = 3.5
    66 = 3.97
    67 = 1.01
    68 = 0.22
    69 = 1.11
    70 = 0.48
    71 = 0.52
    73 = 0.96
    74 = 6.23
    75 = 0.4
    77 = 0.21
    78 = 2.17
    79 = 0.09
    80 = 0.11
    81 = 0.05
    82 = 0.13
    83 = 1.48
    85 = 0.23
    86 = 0.28
    88 = 3.42
}

foreach ($row in $priceUpdates.Keys) {
    $ws.Cells.Item($row, 4).Value = $priceUpdates[$row]
}
